# Update "想去人数" (column F) and "最低票价" (column G) figures for the
# "展览" and "全部类型" sheets to reflect newly scraped data.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Row -> F (想去人数) value
$fUpdates = @{
    2  = 414
    3  = 1402
    4  = 7016
    5  = 530
    6  = 282
    7  = 4624
    8  = 83
    9  = 526
    11 = 932
    12 = 254
    13 = 5437
}

# Row -> G (最低票价) value (only row 5 changes, ticket sold out)
$gUpdates = @{
    5 = "已售罄"
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    foreach ($row in $fUpdates.Keys) {
        $ws.Range("F$row").Value = $fUpdates[$row]
    }

    foreach ($row in $gUpdates.Keys) {
        $ws.Range("G$row").Value = $gUpdates[$row]
    }
}
